$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 790 (shifts existing rows 790-831 down to 791-832)
$ws.Rows("790:790").Insert()

# Seed the new row's formatting/type by copying the row above (789),
# which already carries the same date (2026/02/08) and weekday (日)
# as plain text cells, then overwrite the time value.
$ws.Range("A789:D789").Copy()
$ws.Range("A790:D790").PasteSpecial()
$ws.Range("C790").Value = 21

$excel.CutCopyMode = $false
